# Update "想去人数" (F column) and occasionally "最低票价" (G column)
# figures on the "展览" and "全部类型" worksheets to the newly scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 92
$ws1.Range("F5").Value = 134
$ws1.Range("F6").Value = 9137
$ws1.Range("F7").Value = 830
$ws1.Range("F9").Value = 1173
$ws1.Range("G9").Value = 58
$ws1.Range("F10").Value = 1056
$ws1.Range("F11").Value = 138
$ws1.Range("F12").Value = 52
$ws1.Range("F13").Value = 13
$ws1.Range("F14").Value = 246
$ws1.Range("F15").Value = 352
$ws1.Range("F16").Value = 76
$ws1.Range("F17").Value = 243
$ws1.Range("F18").Value = 1185

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 92
$ws4.Range("F7").Value = 134
$ws4.Range("F8").Value = 9137
$ws4.Range("F9").Value = 830
$ws4.Range("F11").Value = 1173
$ws4.Range("G11").Value = 58
$ws4.Range("F12").Value = 1056
$ws4.Range("F13").Value = 138
$ws4.Range("F14").Value = 52
$ws4.Range("F15").Value = 13
$ws4.Range("F16").Value = 246
$ws4.Range("F17").Value = 352
$ws4.Range("F18").Value = 76
$ws4.Range("F19").Value = 243
$ws4.Range("F20").Value = 1185
